# Append the new INCO claim rows (Caso -507.. -516) reported by the
# 🔄 automatic map update (mapa_interactivo_INCO.html) run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New records, in column order:
# Caso, F.De Reclamo, Direccion, Comuna, OT, Proveedor Asignado, Estado,
# Observaciones, Attachments, Tipo de tarea, Equipo, Tipo de Elemento,
# Coordenada_X, Coordenada_Y, Operacion, Zona
$rows = @(
    @{ Caso="-507"; Fecha="7/14/2025"; Direccion="Tamborini 3291";        Comuna="12"; OT="808194229"; Estado="Pendiente"; Obs="Picada"; Equipo="Sin equipos"; X=-58.473937; Y=-34.557355; Operacion="Saavedra";             Zona="Capital Norte" },
    @{ Caso="-508"; Fecha="7/14/2025"; Direccion="Moldes 2463";           Comuna="12"; OT="808194234"; Estado="Pendiente"; Obs="Picada"; Equipo="Nodo Teco";    X=-58.462281; Y=-34.560321; Operacion="Saavedra";             Zona="Capital Norte" },
    @{ Caso="-509"; Fecha="7/14/2025"; Direccion="Paso 58";               Comuna="3";  OT="808194240"; Estado="Pendiente"; Obs="Picada"; Equipo="Sin equipos"; X=-58.403422; Y=-34.609195; Operacion="Almagro";               Zona="Capital Sur"   },
    @{ Caso="-510"; Fecha="7/14/2025"; Direccion="Larrea 590";            Comuna="3";  OT="808194254"; Estado="Pendiente"; Obs="Picada"; Equipo="Fuente Teco"; X=-58.402353; Y=-34.602205; Operacion="Almagro";               Zona="Capital Sur"   },
    @{ Caso="-512"; Fecha="7/15/2025"; Direccion="Ciudad de la Paz 3742"; Comuna="12"; OT="808240230"; Estado="Pendiente"; Obs="Picada"; Equipo="Sin equipos"; X=-58.470347; Y=-34.547965; Operacion="Saavedra";             Zona="Capital Norte" },
    @{ Caso="-514"; Fecha="7/15/2025"; Direccion="Bilbao 2452";           Comuna="7";  OT="808243829"; Estado="Pendiente"; Obs="Picada"; Equipo="Sin equipos"; X=-58.460594; Y=-34.635581; Operacion="Boedo";                 Zona="Capital Sur"   },
    @{ Caso="-516"; Fecha="7/16/2025"; Direccion="Olazabal 4417";         Comuna="12"; OT="808373646"; Estado="Pendiente"; Obs="Picada"; Equipo="Sin equipos"; X=-58.478941; Y=-34.57242;  Operacion="Colegiales";            Zona="Capital Norte" }
)

$startRow = 30
$endRow = $startRow + $rows.Count - 1

# Columns A, B, D, E hold values that look numeric/date ("-507", "7/14/2025",
# "12", "808194229") but must be stored as plain text, matching the rest of
# the sheet. Pre-formatting the destination columns as Text keeps Excel from
# auto-converting them to numbers / date serials on assignment.
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"
$ws.Range("B$startRow`:B$endRow").NumberFormat = "@"
$ws.Range("D$startRow`:D$endRow").NumberFormat = "@"
$ws.Range("E$startRow`:E$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $rec = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $rec.Caso
    $ws.Cells.Item($r, 2).Value = $rec.Fecha
    $ws.Cells.Item($r, 3).Value = $rec.Direccion
    $ws.Cells.Item($r, 4).Value = $rec.Comuna
    $ws.Cells.Item($r, 5).Value = $rec.OT
    $ws.Cells.Item($r, 6).Value = "INCO"
    $ws.Cells.Item($r, 7).Value = $rec.Estado
    $ws.Cells.Item($r, 8).Value = $rec.Obs
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = "Cambio"
    $ws.Cells.Item($r, 11).Value = $rec.Equipo
    $ws.Cells.Item($r, 12).Value = "Pasante"
    $ws.Cells.Item($r, 13).Value = $rec.X
    $ws.Cells.Item($r, 14).Value = $rec.Y
    $ws.Cells.Item($r, 15).Value = $rec.Operacion
    $ws.Cells.Item($r, 16).Value = $rec.Zona
}
